$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B31").Value = 60
$ws.Range("B79").Value = 17
$ws.Range("B82").Value = 26
$ws.Range("B84").Value = 24
